$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5349307
$ws.Range("I33").Value = 2157.3635
$ws.Range("J33").Value = 15152415
$ws.Range("K33").Value = 2157.3635
$ws.Range("L33").Value = 15152415
$ws.Range("M33").Value = -1928.3635
$ws.Range("N33").Value = -15152873

$ws.Range("H64").Value = 7285.7144
$ws.Range("I64").Value = 10050
$ws.Range("J64").Value = 3600
$ws.Range("K64").Value = 10050
$ws.Range("L64").Value = 3600
$ws.Range("M64").Value = -9802
$ws.Range("N64").Value = -4096

$ws.Range("H67").Value = 7285.7144
$ws.Range("I67").Value = 10050
$ws.Range("J67").Value = 3600
$ws.Range("K67").Value = 10050
$ws.Range("L67").Value = 3600
$ws.Range("M67").Value = -9192
$ws.Range("N67").Value = -5316

$ws.Range("H100").Value = 16667748
$ws.Range("I100").Value = 18519554
$ws.Range("K100").Value = 18519554
$ws.Range("M100").Value = -18519013

$ws.Range("H112").Value = 33615052
$ws.Range("I112").Value = 762.5
$ws.Range("J112").Value = 43957908
$ws.Range("K112").Value = 2287.5
$ws.Range("L112").Value = 131873724
$ws.Range("M112").Value = -1179.5
$ws.Range("N112").Value = -131875940

$ws.Range("H131").Value = 2199.2632
$ws.Range("I131").Value = 1186.25
$ws.Range("K131").Value = 3558.75
$ws.Range("M131").Value = 1481.25

$ws.Range("H132").Value = 4977129.5
$ws.Range("I132").Value = 1437.7288
$ws.Range("J132").Value = 41672860
$ws.Range("K132").Value = 4313.186400000001
$ws.Range("L132").Value = 125018580
$ws.Range("M132").Value = -1783.186400000001
$ws.Range("N132").Value = -125023640

$ws.Range("H138").Value = 4236.9077
$ws.Range("I138").Value = 2319.3076
$ws.Range("J138").Value = 4716.3076
$ws.Range("K138").Value = 6957.9228
$ws.Range("L138").Value = 14148.9228
$ws.Range("M138").Value = -1817.9228
$ws.Range("N138").Value = -24428.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19524.578
$ws.Range("I32").Value = 15616.8125
$ws.Range("J32").Value = 31247.875
$ws.Range("K32").Value = 15616.8125
$ws.Range("L32").Value = 31247.875
$ws.Range("M32").Value = -15329.8125
$ws.Range("N32").Value = -31821.875

$ws.Range("H61").Value = 249051.7
$ws.Range("I61").Value = 6133.7407
$ws.Range("J61").Value = 717536.4
$ws.Range("K61").Value = 6133.7407
$ws.Range("L61").Value = 717536.4
$ws.Range("M61").Value = -5921.7407
$ws.Range("N61").Value = -717960.4

$ws.Range("H74").Value = 9616883
$ws.Range("I74").Value = 1064.4103
$ws.Range("J74").Value = 38464340
$ws.Range("K74").Value = 1064.4103
$ws.Range("L74").Value = 38464340
$ws.Range("M74").Value = -190.4103
$ws.Range("N74").Value = -38466088

$ws.Range("H77").Value = 9616883
$ws.Range("I77").Value = 1064.4103
$ws.Range("J77").Value = 38464340
$ws.Range("K77").Value = 5322.0515
$ws.Range("L77").Value = 192321700
$ws.Range("M77").Value = -954.0514999999996
$ws.Range("N77").Value = -192330436

$ws.Range("H122").Value = 525332.25
$ws.Range("I122").Value = 547559.2
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 1642677.6
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1640227.6
$ws.Range("N122").Value = -13900

$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960

$ws.Range("H132").Value = 2969.1914
$ws.Range("I132").Value = 1854.9032
$ws.Range("K132").Value = 5564.7096
$ws.Range("M132").Value = -3034.7096

$ws.Range("H136").Value = 249051.7
$ws.Range("I136").Value = 6133.7407
$ws.Range("J136").Value = 717536.4
$ws.Range("K136").Value = 18401.2221
$ws.Range("L136").Value = 2152609.2
$ws.Range("M136").Value = -15851.2221
$ws.Range("N136").Value = -2157709.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 523.9091
$ws.Range("I22").Value = 494
$ws.Range("K22").Value = 494
$ws.Range("M22").Value = -321

$ws.Range("H134").Value = 30157.719
$ws.Range("I134").Value = 5019.303
$ws.Range("J134").Value = 168419
$ws.Range("K134").Value = 15057.909
$ws.Range("L134").Value = 505257
$ws.Range("M134").Value = -12522.909
$ws.Range("N134").Value = -510327

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16133937
$ws.Range("I31").Value = 1843.8948
$ws.Range("J31").Value = 41676420
$ws.Range("K31").Value = 1843.8948
$ws.Range("L31").Value = 41676420
$ws.Range("M31").Value = -1548.8948
$ws.Range("N31").Value = -41677010

$ws.Range("H34").Value = 16133937
$ws.Range("I34").Value = 1843.8948
$ws.Range("J34").Value = 41676420
$ws.Range("K34").Value = 1843.8948
$ws.Range("L34").Value = 41676420
$ws.Range("M34").Value = -1641.8948
$ws.Range("N34").Value = -41676824

$ws.Range("H86").Value = 2175.3845
$ws.Range("I86").Value = 2185.7144
$ws.Range("K86").Value = 2185.7144
$ws.Range("M86").Value = -1062.7144

$ws.Range("H89").Value = 2175.3845
$ws.Range("I89").Value = 2185.7144
$ws.Range("K89").Value = 10928.572
$ws.Range("M89").Value = -5312.572

$ws.Range("H132").Value = 5264722
$ws.Range("J132").Value = 3925.25
$ws.Range("L132").Value = 11775.75
$ws.Range("N132").Value = -16835.75

$ws.Range("H134").Value = 12265288
$ws.Range("I134").Value = 13337017
$ws.Range("J134").Value = 3334208.2
$ws.Range("K134").Value = 40011051
$ws.Range("L134").Value = 10002624.6
$ws.Range("M134").Value = -40008516
$ws.Range("N134").Value = -10007694.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10614.883
$ws.Range("I3").Value = 16465
$ws.Range("J3").Value = 5414.778
$ws.Range("K3").Value = 49395
$ws.Range("L3").Value = 16244.334
$ws.Range("M3").Value = -49283
$ws.Range("N3").Value = -16468.334

$ws.Range("H107").Value = 501.1143
$ws.Range("I107").Value = 339.9524
$ws.Range("J107").Value = 742.8570999999999
$ws.Range("K107").Value = 1019.8572
$ws.Range("L107").Value = 2228.5713
$ws.Range("M107").Value = 900.1428
$ws.Range("N107").Value = -6068.5713

$ws.Range("H114").Value = 5392.6665
$ws.Range("I114").Value = 345.66666
$ws.Range("J114").Value = 7916.1665
$ws.Range("K114").Value = 1036.99998
$ws.Range("L114").Value = 23748.4995
$ws.Range("M114").Value = 2217.00002
$ws.Range("N114").Value = -30256.4995

$ws.Range("H122").Value = 8593.786
$ws.Range("I122").Value = 1147.6666
$ws.Range("J122").Value = 14178.375
$ws.Range("K122").Value = 10328.9994
$ws.Range("L122").Value = 127605.375
$ws.Range("M122").Value = -7878.999400000001
$ws.Range("N122").Value = -132505.375

$ws.Range("H131").Value = 1887977.8
$ws.Range("J131").Value = 1289.4762
$ws.Range("L131").Value = 3868.4286
$ws.Range("N131").Value = -13948.4286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6668829
$ws.Range("I132").Value = 8334986.5
$ws.Range("K132").Value = 25004959.5
$ws.Range("M132").Value = -25002429.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2864

$ws.Range("H68").Value = 41668704
$ws.Range("I68").Value = 1898.6
$ws.Range("J68").Value = 111113380
$ws.Range("K68").Value = 1898.6
$ws.Range("L68").Value = 111113380
$ws.Range("M68").Value = -1149.6
$ws.Range("N68").Value = -111114878

$ws.Range("H71").Value = 41668704
$ws.Range("I71").Value = 1898.6
$ws.Range("J71").Value = 111113380
$ws.Range("K71").Value = 9493
$ws.Range("L71").Value = 555566900
$ws.Range("M71").Value = -5749
$ws.Range("N71").Value = -555574388

$ws.Range("H122").Value = 3881245.8
$ws.Range("I122").Value = 4468229
$ws.Range("J122").Value = 2002899
$ws.Range("K122").Value = 13404687
$ws.Range("L122").Value = 6008697
$ws.Range("M122").Value = -13402237
$ws.Range("N122").Value = -6013597

$ws.Range("H132").Value = 7136.032
$ws.Range("I132").Value = 7214.759
$ws.Range("J132").Value = 5994.5
$ws.Range("K132").Value = 21644.277
$ws.Range("L132").Value = 17983.5
$ws.Range("M132").Value = -19114.277
$ws.Range("N132").Value = -23043.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3225
$ws.Range("I96").Value = 3225
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 3225
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -1852
$ws.Range("N96").ClearContents()
